$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." footer line.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # The paragraph immediately before it is the blank spacer paragraph that
    # was introduced together with the footer block, and the paragraph right
    # after it is the "(c) 2020 ..." copyright line - both need to go too.
    $startPara = $d.Paragraphs.Item($target - 1)
    $endPara = $d.Paragraphs.Item($target + 1)

    $rangeToDelete = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rangeToDelete.Delete()
}
